$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos table (columns B..E) with the latest scraped values.
# Column D holds numeric-looking strings (e.g. "44.124.20", "0.0000100") that must stay
# plain text exactly as scraped, so we force a text number format while assigning them,
# then restore the cell style back to Normal so no extra formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.124.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.276.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.621.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.11%  "

$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.281.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.068.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.04%  "

$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0888"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "161.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.39%  "

$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.53%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.108"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.91%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +37.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0328"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.798.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.06%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
